$d = $word.ActiveDocument

# wdCharacter unit constant used by Range.MoveEnd/MoveStart
$wdCharacter = 1

# --- Paragraph 1 ("The Eleventh Virgin" / "Part I, Chapter IV =====...") ---
# This whole paragraph (title + line break + chapter marker) is removed.
$titlePara = $d.Paragraphs(1)
$null = $titlePara.Range.Delete()

# --- Former paragraph 2, now paragraph 1 ("By Dorothy Day", bold) ---
# Replace its run with an unformatted "% Dorothy Day" (pandoc-style title block).
$byLine = $d.Paragraphs(1)
$byLineText = $byLine.Range
$null = $byLineText.MoveEnd($wdCharacter, -1)   # exclude the trailing paragraph mark
$null = $byLineText.Delete()                     # drop the old "By Dorothy Day" (bold) text
$null = $byLineText.InsertAfter("% Dorothy Day") # insert plain, unformatted text
